$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035536927939108
$ws.Range("D2").Value = 1.042979261107844
$ws.Range("E2").Value = 1.039183709549561
$ws.Range("F2").Value = 1.050224574873085
$ws.Range("I2").Value = 1.0361557452688
$ws.Range("J2").Value = 1.040650162953251
$ws.Range("K2").Value = 1.045754420962579
$ws.Range("L2").Value = 1.041969620700722
$ws.Range("M2").Value = 1.052979445582076
$ws.Range("N2").Value = 1.042128005548887

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036709142268484
$ws.Range("D3").Value = 1.043916026199733
$ws.Range("E3").Value = 1.040304627735821
$ws.Range("F3").Value = 1.051358372961392
$ws.Range("I3").Value = 1.036422543245636
$ws.Range("J3").Value = 1.041464843160206
$ws.Range("K3").Value = 1.046501798447456
$ws.Range("L3").Value = 1.042899885740154
$ws.Range("M3").Value = 1.053924818824181
$ws.Range("N3").Value = 1.042943842695181

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03746737792507
$ws.Range("D4").Value = 1.044521555393096
$ws.Range("E4").Value = 1.041030031262062
$ws.Range("F4").Value = 1.052091873222573
$ws.Range("I4").Value = 1.036593173753861
$ws.Range("J4").Value = 1.041991217203077
$ws.Range("K4").Value = 1.04698415925924
$ws.Range("L4").Value = 1.04350134172761
$ws.Range("M4").Value = 1.054535812974521
$ws.Range("N4").Value = 1.043470964249568

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037786078238184
$ws.Range("D5").Value = 1.044775971804059
$ws.Range("E5").Value = 1.041335014489248
$ws.Range("F5").Value = 1.052400203384752
$ws.Range("I5").Value = 1.036664426773454
$ws.Range("J5").Value = 1.042212319221341
$ws.Range("K5").Value = 1.047186646826477
$ws.Range("L5").Value = 1.0437540779041
$ws.Range("M5").Value = 1.054792501943669
$ws.Range("N5").Value = 1.043692380258059

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037839585835766
$ws.Range("D6").Value = 1.044818680770228
$ws.Range("E6").Value = 1.041386223925898
$ws.Range("F6").Value = 1.052451971445235
$ws.Range("I6").Value = 1.036676362321735
$ws.Range("J6").Value = 1.042249432359438
$ws.Range("K6").Value = 1.047220627983367
$ws.Range("L6").Value = 1.043796506646749
$ws.Range("M6").Value = 1.054835591028915
$ws.Range("N6").Value = 1.04372954610107

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037471636655735
$ws.Range("D7").Value = 1.044524955498842
$ws.Range("E7").Value = 1.041034106369012
$ws.Range("F7").Value = 1.05209599327526
$ws.Range("I7").Value = 1.036594127726338
$ws.Range("J7").Value = 1.041994172307541
$ws.Range("K7").Value = 1.046986866074092
$ws.Range("L7").Value = 1.043504719255232
$ws.Range("M7").Value = 1.054539243543864
$ws.Range("N7").Value = 1.04347392355062

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.035933139165111
$ws.Range("D8").Value = 1.043295974017902
$ws.Range("E8").Value = 1.039562510110502
$ws.Range("F8").Value = 1.050607777361988
$ws.Range("I8").Value = 1.036246326236746
$ws.Range("J8").Value = 1.04092564959731
$ws.Range("K8").Value = 1.04600725824179
$ws.Range("L8").Value = 1.042284109612358
$ws.Range("M8").Value = 1.053299089566449
$ws.Range("N8").Value = 1.042403883415578

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.033220006089745
$ws.Range("D9").Value = 1.041125577239744
$ws.Range("E9").Value = 1.036970049574725
$ws.Range("F9").Value = 1.047984216490849
$ws.Range("I9").Value = 1.035618089559687
$ws.Range("J9").Value = 1.039036779122885
$ws.Range("K9").Value = 1.044271524884969
$ws.Range("L9").Value = 1.040129466169788
$ws.Range("M9").Value = 1.051108189503041
$ws.Range("N9").Value = 1.040512330528424

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.031409733653298
$ws.Range("D10").Value = 1.039675401120386
$ws.Range("E10").Value = 1.035242135519953
$ws.Range("F10").Value = 1.046234354201024
$ws.Range("I10").Value = 1.035188924279496
$ws.Range("J10").Value = 1.037773449358778
$ws.Range("K10").Value = 1.043107919112659
$ws.Range("L10").Value = 1.038690444041887
$ws.Range("M10").Value = 1.049643780776949
$ws.Range("N10").Value = 1.039247206691146

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030625482199473
$ws.Range("D11").Value = 1.039046680319312
$ws.Range("E11").Value = 1.034494007580294
$ws.Range("F11").Value = 1.04547643537124
$ws.Range("I11").Value = 1.035000635464629
$ws.Range("J11").Value = 1.037225433087174
$ws.Range("K11").Value = 1.042602525949046
$ws.Range("L11").Value = 1.038066703093607
$ws.Range("M11").Value = 1.049008760088854
$ws.Range("N11").Value = 1.038698412173578

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030334115318575
$ws.Range("D12").Value = 1.038813026596425
$ws.Range("E12").Value = 1.034216128680007
$ws.Range("F12").Value = 1.045194876632896
$ws.Range("I12").Value = 1.034930327083526
$ws.Range("J12").Value = 1.037021726134943
$ws.Range("K12").Value = 1.042414567525448
$ws.Range("L12").Value = 1.037834921288166
$ws.Range("M12").Value = 1.048772745434297
$ws.Range("N12").Value = 1.038494415934109

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.030396617297286
$ws.Range("D13").Value = 1.038863151504579
$ws.Range("E13").Value = 1.034275734257404
$ws.Range("F13").Value = 1.045255273495638
$ws.Range("I13").Value = 1.034945425197443
$ws.Range("J13").Value = 1.037065428742794
$ws.Range("K13").Value = 1.042454895799184
$ws.Range("L13").Value = 1.037884643666007
$ws.Range("M13").Value = 1.048823377719334
$ws.Range("N13").Value = 1.038538180604676

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.030601398978468
$ws.Range("D14").Value = 1.039027368858845
$ws.Range("E14").Value = 1.034471037845812
$ws.Range("F14").Value = 1.045453162338281
$ws.Range("I14").Value = 1.034994831292197
$ws.Range("J14").Value = 1.037208597659492
$ws.Range("K14").Value = 1.04258699400347
$ws.Range("L14").Value = 1.038047545916132
$ws.Range("M14").Value = 1.048989253900455
$ws.Range("N14").Value = 1.038681552837657

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030727563587116
$ws.Range("D15").Value = 1.039128532818553
$ws.Range("E15").Value = 1.034591371999018
$ws.Range("F15").Value = 1.045575083647618
$ws.Range("I15").Value = 1.035025223041176
$ws.Range("J15").Value = 1.037296788934835
$ws.Range("K15").Value = 1.042668353178203
$ws.Range("L15").Value = 1.038147902524534
$ws.Range("M15").Value = 1.049091437143912
$ws.Range("N15").Value = 1.038769869354727

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031461773406977
$ws.Range("D16").Value = 1.039717110635126
$ws.Range("E16").Value = 1.035291787665861
$ws.Range("F16").Value = 1.046284650092553
$ws.Range("I16").Value = 1.035201368584852
$ws.Range("J16").Value = 1.037809798518834
$ws.Range("K16").Value = 1.043141427781356
$ws.Range("L16").Value = 1.03873182619078
$ws.Range("M16").Value = 1.04968590548258
$ws.Range("N16").Value = 1.039283607471178

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.031922217489819
$ws.Range("D17").Value = 1.040086099165813
$ws.Range("E17").Value = 1.035731157712447
$ws.Range("F17").Value = 1.046729683354633
$ws.Range("I17").Value = 1.035311201878159
$ws.Range("J17").Value = 1.03813133104832
$ws.Range("K17").Value = 1.043437760925015
$ws.Range("L17").Value = 1.039097935218268
$ws.Range("M17").Value = 1.050058552085993
$ws.Range("N17").Value = 1.039605596613736

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.032190749099567
$ws.Range("D18").Value = 1.04030124825499
$ws.Range("E18").Value = 1.035987441701302
$ws.Range("F18").Value = 1.046989242889485
$ws.Range("I18").Value = 1.035375028697566
$ws.Range("J18").Value = 1.03831878046167
$ws.Range("K18").Value = 1.043610458103836
$ws.Range("L18").Value = 1.039311419249321
$ws.Range("M18").Value = 1.050275821660072
$ws.Range("N18").Value = 1.039793312226754

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03228230512087
$ws.Range("D19").Value = 1.040374595725056
$ws.Range("E19").Value = 1.036074829127787
$ws.Range("F19").Value = 1.047077742489739
$ws.Range("I19").Value = 1.035396751790285
$ws.Range("J19").Value = 1.038382679776657
$ws.Range("K19").Value = 1.043669318145082
$ws.Range("L19").Value = 1.039384201449678
$ws.Range("M19").Value = 1.050349889929063
$ws.Range("N19").Value = 1.039857302286098

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.031872820105155
$ws.Range("D20").Value = 1.040046518022469
$ws.Range("E20").Value = 1.035684016763618
$ws.Range("F20").Value = 1.046681937662171
$ws.Range("I20").Value = 1.03529944233006
$ws.Range("J20").Value = 1.038096843492611
$ws.Range("K20").Value = 1.043405982594039
$ws.Range("L20").Value = 1.039058661520487
$ws.Range("M20").Value = 1.050018579834133
$ws.Range("N20").Value = 1.039571060081743

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.030541097577828
$ws.Range("D21").Value = 1.038979014228347
$ws.Range("E21").Value = 1.0344135255664
$ws.Range("F21").Value = 1.045394889954145
$ws.Range("I21").Value = 1.034980292632298
$ws.Range("J21").Value = 1.037166442112849
$ws.Range("K21").Value = 1.042548100811258
$ws.Range("L21").Value = 1.03799957794177
$ws.Range("M21").Value = 1.048940411365231
$ws.Range("N21").Value = 1.038639337425303

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029703436933581
$ws.Range("D22").Value = 1.038307144297954
$ws.Range("E22").Value = 1.033614769417213
$ws.Range("F22").Value = 1.044585474942416
$ws.Range("I22").Value = 1.034777492384059
$ws.Range("J22").Value = 1.036580597408721
$ws.Range("K22").Value = 1.042007368542669
$ws.Range("L22").Value = 1.037333130270031
$ws.Range("M22").Value = 1.048261715227616
$ws.Range("N22").Value = 1.038052660754498

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.030147530949196
$ws.Range("D23").Value = 1.03866338076828
$ws.Range("E23").Value = 1.034038200706311
$ws.Range("F23").Value = 1.045014580285376
$ws.Range("I23").Value = 1.034885203419226
$ws.Range("J23").Value = 1.03689124709559
$ws.Range("K23").Value = 1.042294149023258
$ws.Range("L23").Value = 1.037686480276372
$ws.Range("M23").Value = 1.048621581952712
$ws.Range("N23").Value = 1.038363751599555

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03189514077269
$ws.Range("D24").Value = 1.040064403270882
$ws.Range("E24").Value = 1.03570531770755
$ws.Range("F24").Value = 1.046703511949346
$ws.Range("I24").Value = 1.035304756696038
$ws.Range("J24").Value = 1.038112427227915
$ws.Range("K24").Value = 1.043420342313757
$ws.Range("L24").Value = 1.0390764078026
$ws.Range("M24").Value = 1.050036641848402
$ws.Range("N24").Value = 1.039586665947738

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.033921675872909
$ws.Range("D25").Value = 1.04168724616891
$ws.Range("E25").Value = 1.037640188755406
$ws.Range("F25").Value = 1.048662610323024
$ws.Range("I25").Value = 1.035782325709454
$ws.Range("J25").Value = 1.039525813280932
$ws.Range("K25").Value = 1.044721387732643
$ws.Range("L25").Value = 1.040686946056841
$ws.Range("M25").Value = 1.051675257280717
$ws.Range("N25").Value = 1.041002059171068

